$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting rows 13:17 down to 14:18
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new weekly data entry
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 44596
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 100112001
$ws.Cells.Item(13, 7).Value = "Berenjena"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 12000
$ws.Cells.Item(13, 12).Value = 13000
$ws.Cells.Item(13, 13).Value = 12500
$ws.Cells.Item(13, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 208
$ws.Cells.Item(13, 17).Value = 60
$ws.Cells.Item(13, 18).Value = "Hortaliza"
